$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.147.59"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").Value = "3.293.19"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'586.17"
$ws.Range("E5").Value = "  +2.97%  "

$ws.Range("D6").Value = "'180.54"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").Value = "'0.656"
$ws.Range("E7").Value = "  +9.81%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -2.00%  "

$ws.Range("D10").Value = "'6.75"
$ws.Range("E10").Value = "  +1.95%  "

$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("D12").Value = "3.864.99"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("E13").Value = "  -5.15%  "

$ws.Range("D14").Value = "66.173.04"
$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("D15").Value = "'26.46"
$ws.Range("E15").Value = "  -1.71%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000163"
$ws.Range("E16").Value = "  -1.45%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.294.11"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("D18").Value = "'436.04"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("D19").Value = "'13.26"
$ws.Range("E19").Value = "  -2.35%  "

$ws.Range("D20").Value = "'5.49"
$ws.Range("E20").Value = "  -2.80%  "

$ws.Range("D21").Value = "'7.44"
$ws.Range("E21").Value = "  -2.39%  "

$ws.Range("D22").Value = "'72.39"
$ws.Range("E22").Value = "  -1.22%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "'5.68"
$ws.Range("E24").Value = "  +0.41%  "

$ws.Range("D25").Value = "3.427.94"
$ws.Range("E25").Value = "  -0.61%  "

$ws.Range("D26").Value = "'0.511"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "'0.198"
$ws.Range("E27").Value = "  +3.92%  "

$ws.Range("E28").Value = "  -2.97%  "

$ws.Range("D29").Value = "'8.86"
$ws.Range("E29").Value = "  -0.84%  "

$ws.Range("E30").Value = "  +0.27%  "

$ws.Range("E31").Value = "  +1.10%  "

$ws.Range("D32").Value = "'22.34"
$ws.Range("E32").Value = "  -1.42%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").Value = "'5.21"
$ws.Range("E34").Value = "  -0.92%  "

$ws.Range("D35").Value = "'6.62"
$ws.Range("E35").Value = "  -1.63%  "

$ws.Range("D36").Value = "'1.19"
$ws.Range("E36").Value = "  -1.65%  "

$ws.Range("D37").Value = "'158.09"
$ws.Range("E37").Value = "  -0.74%  "

$ws.Range("E38").Value = "  -4.74%  "

$ws.Range("D39").Value = "'26.55"
$ws.Range("E39").Value = "  -2.12%  "

$ws.Range("E40").Value = "  -2.92%  "

$ws.Range("D41").Value = "2.800.49"
$ws.Range("E41").Value = "  +0.95%  "

$ws.Range("D42").Value = "'0.774"
$ws.Range("E42").Value = "  -1.03%  "

$ws.Range("D43").Value = "'4.35"
$ws.Range("E43").Value = "  -1.69%  "

$ws.Range("D44").Value = "'40.24"
$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("D45").Value = "'6.10"
$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("D46").Value = "'0.0661"
$ws.Range("E46").Value = "  -1.50%  "

$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("D48").Value = "'320.69"
$ws.Range("E48").Value = "  +0.90%  "

$ws.Range("D49").Value = "'23.26"
$ws.Range("E49").Value = "  -3.37%  "

$ws.Range("E50").Value = "  -0.99%  "

$ws.Range("E51").Value = "  +7.20%  "
